# Insert a new weekly record at row 107 (shifting the existing rows
# 107..230 down to 108..231, like Excel's "Insert Sheet Rows").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

$ws.Cells.Item(107, 1).Value = 6
$ws.Cells.Item(107, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(107, 3).Value = "Metropolitana"
$ws.Cells.Item(107, 4).Value = 44803
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 100112029
$ws.Cells.Item(107, 7).Value = "Orégano"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 47
$ws.Cells.Item(107, 11).Value = 15000
$ws.Cells.Item(107, 12).Value = 16000
$ws.Cells.Item(107, 13).Value = 15447
$ws.Cells.Item(107, 14).Value = "$/docena de atados"
$ws.Cells.Item(107, 15).Value = "Región Metropolitana"
$ws.Cells.Item(107, 16).Value = 5149
$ws.Cells.Item(107, 17).Value = 3
$ws.Cells.Item(107, 18).Value = "Hortaliza"
